# Ajout accent police bebas
# Fill in the Iteration #1 worklog (dates/tasks/hours for the first week),
# the end-of-iteration hour total, the self-evaluation comment, and make
# "Iteration #1" the active/selected sheet (scrolled down to the totals).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #1")

# --- Daily worklog rows 14-18: date / task description / hours worked ---
$rows = @(
    @{ Row = 14; Date = 43179; Task = "Création de l'app, création du logo, développement de ma vue du main menu"; Hours = 2.5 },
    @{ Row = 15; Date = 43180; Task = "Recherche et implémentation du expandableListView avec pager adapter (données hard code pour l'instant)"; Hours = 4.5 },
    @{ Row = 16; Date = 43181; Task = "Création du modèle relationnel de ma BD + fin implémentation listView"; Hours = 3 },
    @{ Row = 17; Date = 43186; Task = "Début de la création de mon modèle SQLITE dans mon app android"; Hours = 2 },
    @{ Row = 18; Date = 43187; Task = "Fin modèle SQLITE + début création BD avec sqlite data browser"; Hours = 3 }
)

# Format the first date cell (center aligned, short-date number format), then
# fan that exact formatting out to the other date cells via copy/paste-format
# so they all share a single cell style instead of one each.
$firstDateCell = $ws.Cells.Item(14, 1)
$firstDateCell.Value = 43179
$firstDateCell.HorizontalAlignment = -4108   # xlCenter
$firstDateCell.NumberFormat = "mm-dd-yy"     # builtin short date format

$firstDateCell.Copy()
$ws.Range("A15:A18").PasteSpecial(-4122)     # xlPasteFormats

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Date
    $ws.Cells.Item($r.Row, 2).Value = $r.Task
    $ws.Cells.Item($r.Row, 3).Value = $r.Hours
}

# The remaining still-empty date cells (rows 19-36) just get re-centered,
# matching the formatting used for the filled-in rows above.
$ws.Range("A19:A36").HorizontalAlignment = -4108

# --- End of iteration summary ---
# Total hours spent this iteration
$ws.Cells.Item(40, 2).Value = 7

# Self-evaluation comment
$ws.Cells.Item(42, 2).Value = "J'ai bien enclenché le projet, par contre, je pourrais mettre plus d'heure pour la prochaine itération."

# Recalculate so the SUM(C14:C36) total (row 37) reflects the new hours
$excel.Calculate()

# --- Make "Iteration #1" the active sheet, scrolled/selected near the bottom ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B42:B47").Select()
